$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper behaviour note: assigning Range.Text with a value that differs from
# the range's current text causes adjoining runs that share the same
# formatting to be merged into a single run (while a no-op assignment, i.e.
# setting the text to the value it already has, is treated as a no-op and
# does not trigger a merge). So wherever the XML only collapses several runs
# that already contain the same characters into one run, we first overwrite
# the target run with a placeholder string and then restore the original
# text; this forces Word to re-author that span as one run without altering
# the visible text.
# ---------------------------------------------------------------------------

# --- Change 1: paragraph "+ 1_Anh he thong ..." -> merge "+ ", "1_Ảnh hệ
# thống" and the trailing sentence into a single run. ---
$p2 = $d.Paragraphs(2)
$r1 = $p2.Range.Duplicate
$r1.Find.Execute("1_Ảnh hệ thống") | Out-Null
$r1.Text = "@@@@@@@@@@@@@@"
$r1b = $p2.Range.Duplicate
$r1b.Find.Execute("@@@@@@@@@@@@@@") | Out-Null
$r1b.Text = "1_Ảnh hệ thống"

# --- Change 2: paragraph "+ 2_Tu lieu he thong frontend ..." -> merge the
# ", neu dinh dang ... PDF co ten "", "Tư liệu hệ thống" and """. " runs
# into a single run (the quoted "Tư liệu hệ thống" here is the *second*
# occurrence of that phrase in the paragraph -- the first is part of the
# heading "2_Tư liệu hệ thống frontend"). ---
$p3 = $d.Paragraphs(3)
$r2 = $p3.Range.Duplicate
$r2.Find.Execute("Tư liệu hệ thống") | Out-Null
$r2.Find.Execute("Tư liệu hệ thống") | Out-Null
$r2.Text = "@@@@@@@@@@@@@@@@"
$r2b = $p3.Range.Duplicate
$r2b.Find.Execute("@@@@@@@@@@@@@@@@") | Out-Null
$r2b.Text = "Tư liệu hệ thống"

# --- Change 3: insert a brand-new paragraph "+ 3_Tai khoan 4 cap: Chu tai
# khoan va mat khau 4 cap cua app" right before the "+ Nhung file, anh khong
# duoc danh so ..." paragraph, then re-merge that paragraph's first two runs
# (the intro text + the quoted "Tư liệu hệ thống") into one run. ---
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs(4)
$newPara.Range.Text = "+ 3_Tài khoản 4 cấp: Chứ tài khoản và mật khẩu 4 cấp của app"

$p5 = $d.Paragraphs(5)
$r3 = $p5.Range.Duplicate
$r3.Find.Execute("Tư liệu hệ thống") | Out-Null
$r3.Text = "@@@@@@@@@@@@@@@@"
$r3b = $p5.Range.Duplicate
$r3b.Find.Execute("@@@@@@@@@@@@@@@@") | Out-Null
$r3b.Text = "Tư liệu hệ thống"

# --- Change 4: hyperlink text "Xm" + "i" + "nd" -> single run "Xmind". ---
$r4 = $d.Content
$r4.Find.Execute("Xmind", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$midStart = $r4.Start + 2
$midEnd = $r4.Start + 3
$r4mid = $d.Range($midStart, $midEnd)
$r4mid.Text = "@"
$r4mid2 = $d.Range($midStart, $midStart + 1)
$r4mid2.Text = "i"
